# "Generate Report for Handoff"
#
# Marks the localization-status report as ready for a fresh handoff:
#   - status cells flip from "Handed back: in sync with en-US" to "Ready for handoff"
#   - the handoff-generation timestamps are bumped to the new run time
#   - the (now shorter) status column is narrowed to fit the new text
#
# Column E/F (Overview) and column C (zh-cn / de-de) were narrowed from the
# old "Handed back..." width down to a tighter width. Excel's ColumnWidth
# COM setter only accepts values on a 1/6-character pixel grid (it rounds
# internally), so the nearest representable width is used here.
$narrowWidth = 16.333333333333332

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-27 22:58:27"
$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-27 22:58:23"
$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-27 22:58:27"
$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
